$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 0. The old "params" header text becomes unused once the column is repurposed as "query" below;
#     clear it now (while it is still the sole reference) so the shared-string table drops the entry.
$ws.Range("C1").ClearContents()

# --- 1. Insert a new column C ("body"), shifting the old C (now "query") and D (now "response") to the right.
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).ColumnWidth = 15.71

# --- 2. Insert the new POST row (after "/users" GET) and PUT row (after "/users/1" GET).
# Final layout (1-indexed rows):
#  1 header
#  2 GET /hello
#  3 GET /users
#  4 POST /users        <-- new
#  5 GET /users/1
#  6 PUT /users/1        <-- new
#  7 GET /users/1/books/1
#  8 GET /entitlements (krm)
#  9 GET /entitlements (krr)
# 10 GET /wallets/sample
# 11 GET /wallets
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(6).Insert()

# New POST/PUT rows need the bordered cell style in columns A and B (copy from an existing data cell);
# C/E are intentionally left unstyled to match the target workbook.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Row 4: POST /users (new) -- fill before row 6 / before the header rename so new shared
#     strings are appended to the table in the same order the author would have typed them.
$ws.Range("A4").Value2 = "POST"
$ws.Range("B4").Value2 = "/users"
$ws.Range("C4").Value2 = "{`n  ""name"": ""Pol"",`n  ""lastName"": ""Puig""`n}"
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value2 = "{`n  ""id"": ""1""`n}"
$ws.Rows.Item(4).AutoFit()

# --- 4. Row 6: PUT /users/1 (new)
$ws.Range("A6").Value2 = "PUT"
$ws.Range("B6").Value2 = "/users/1"
$ws.Range("C6").Value2 = "{`n  ""id"": ""1"",`n  ""name"": ""Pol"",`n  ""lastName"": ""New Last Name""`n}"
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value2 = "{`n  ""id"": ""1"",`n  ""name"": ""Pol"",`n  ""lastName"": ""New Last Name""`n}"
$ws.Rows.Item(6).AutoFit()

# --- 5. Header row (renamed last, matching the new shared-string append order)
$ws.Range("A1").Value2 = "method"
$ws.Range("B1").Value2 = "url"
$ws.Range("C1").Value2 = "body"
$ws.Range("D1").Value2 = "query"
$ws.Range("E1").Value2 = "response"

# --- 6. Remaining (pre-existing) rows -- only column C is new/blank for these, D/E keep their data.
$ws.Range("C2").Value2 = ""
$ws.Range("C3").Value2 = ""
$ws.Range("C5").Value2 = ""
$ws.Range("C7").Value2 = ""
$ws.Range("C8").Value2 = ""
$ws.Range("C9").Value2 = ""
$ws.Range("C10").Value2 = ""
$ws.Range("C11").Value2 = ""

# --- 7. Restore the selection/active cell the author left the sheet on.
$ws.Range("D24").Select()
